$d = $word.ActiveDocument

# --- Fill in the empty footnotes (bug: footnotes with no real content) ---

# Footnote id="21" (1st footnote) : "।" -> "སཱ་དྷ་ན་ཀྲ་མ། ཞེས་པར་མ་གཞན་ནང་མེད།"
$d.Footnotes.Item(1).Range.Text = " སཱ་དྷ་ན་ཀྲ་མ། ཞེས་པར་མ་གཞན་ནང་མེད།"

# Footnote id="27" (7th footnote) : drop the stray trailing "a"
$d.Footnotes.Item(7).Range.Text = " གི། པེ་ཅིན། གིས། སྣར་ཐང་།"

# Footnote id="41" (21st footnote) : "।" -> "ཁཱ་ཧི། ཞེས་པར་མ་གཞན་ནང་མེད།"
$d.Footnotes.Item(21).Range.Text = " ཁཱ་ཧི། ཞེས་པར་མ་གཞན་ནང་མེད།"

# Footnote id="42" (22nd footnote) : "।" -> "ཚར་གཅོད་བྱང་ཆུབ་སེམས་དཔའ་སྟེ།_། ཞེས་པར་མ་གཞན་ནང་མེད།"
$d.Footnotes.Item(22).Range.Text = " ཚར་གཅོད་བྱང་ཆུབ་སེམས་དཔའ་སྟེ།_། ཞེས་པར་མ་གཞན་ནང་མེད།"

# --- Remove the trailing, completely empty footnote (id="71") and its reference ---
$d.Footnotes.Item($d.Footnotes.Count).Delete()
